# Applies the 2017-01-31 EIA Table A.2.B update (EPM_2016_11 / November 2016 run).
# Source: www.eia.gov/electricity/monthly/xls/table_a_2_b.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata: refresh the report title (October -> November 2016) ---
$ws.Range("A2").Value = "Electric Utilities by Census Division and State, Year-to-Date through November 2016"

# --- Refreshed Relative Standard Error figures for the November 2016 data run ---
# Maps each changed cell (column B-H, rows 4-64) to its new value.
$cellUpdates = [ordered]@{
    "C4" = 27
    "E4" = 0
    "H4" = 41
    "C5" = 33
    "H5" = 327
    "C6" = 405
    "C7" = 64
    "E7" = 0
    "H7" = 112
    "C8" = 234
    "H8" = 32
    "C9" = 32
    "C10" = 561
    "H10" = 75
    "C11" = 131
    "E11" = 11
    "H11" = 1
    "C12" = 719
    "E12" = 270
    "C13" = 133
    "E13" = 11
    "C14" = 176
    "E14" = 0
    "H14" = 197
    "C15" = 5
    "F15" = 36
    "H15" = 16
    "C16" = 51
    "E16" = 0
    "H16" = 240
    "C17" = 5
    "E17" = 1
    "F17" = 941
    "H17" = 15
    "C18" = 14
    "H18" = 31
    "C19" = 3
    "E19" = 1
    "H19" = 6
    "C20" = 23
    "E20" = 2
    "H20" = 25
    "C21" = 9
    "H21" = 8
    "C22" = 20
    "E22" = 15
    "H22" = 37
    "C23" = 26
    "E23" = 13
    "C24" = 56
    "E24" = 3
    "H24" = 48
    "E25" = 14
    "H25" = 13
    "C26" = 128
    "E26" = 8
    "H26" = 28
    "C27" = 10
    "E27" = 43
    "C28" = 532
    "E28" = 14
    "H28" = 0.44
    "B29" = 0
    "C29" = 9
    "E29" = 0.31
    "H29" = 11
    "C30" = 683
    "E30" = 141
    "C31" = 8
    "E31" = 1
    "H31" = 92
    "B32" = 0
    "C32" = 82
    "E32" = 0
    "H32" = 15
    "C33" = 47
    "C34" = 35
    "H34" = 17
    "C35" = 78
    "E35" = 0
    "H35" = 28
    "C36" = 30
    "E36" = 0.07
    "H36" = 38
    "H37" = 69
    "B38" = 0.28
    "C38" = 5
    "E38" = 2
    "H38" = 7
    "E39" = 6
    "H39" = 11
    "B40" = 1
    "E40" = 0
    "H40" = 7
    "C41" = 437
    "C42" = 0.26
    "H42" = 12
    "H43" = 14
    "E44" = 3
    "H44" = 16
    "C46" = 6
    "H46" = 30
    "E47" = 2
    "H47" = 41
    "C48" = 12
    "H48" = 4
    "C49" = 9
    "C50" = 87
    "E50" = 0
    "H50" = 36
    "C51" = 475
    "E51" = 30
    "H51" = 10
    "B52" = 162
    "C52" = 1551
    "E52" = 62
    "H52" = 4
    "C54" = 56
    "E54" = 7
    "H54" = 116
    "C55" = 4
    "E55" = 7
    "H55" = 42
    "E56" = 140
    "H56" = 27
    "C57" = 32
    "E57" = 3
    "C58" = 12
    "E58" = 4
    "H58" = 6
    "E59" = 0.2
    "C60" = 324
    "E60" = 5
    "C61" = 5
    "E61" = 13
    "H61" = 20
    "C62" = 6
    "E62" = 13
    "H62" = 20
    "C63" = 6
    "H63" = 244
    "B64" = 0.22
    "C64" = 4
    "E64" = 0.42
    "F64" = 36
}

foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
